$d = $word.ActiveDocument

# Locate the three target paragraphs by their distinctive leading text so we
# do not rely on a brittle fixed paragraph index.
$paraStack = $null
$paraProc = $null
$paraStatic = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Stack manipulation operations will be implemented at the syntactic sugar level*") {
        $paraStack = $p
    } elseif ($t -like "Procedures/functions with arguments will be implemented at the syntactic sugar level*") {
        $paraProc = $p
    } elseif ($t -like "Static type system will be implemented at the syntactic sugar level*") {
        $paraStatic = $p
    }
}

# Note: MatchWholeWord (4th argument) must be $false for search strings that
# start/end with a space or punctuation -- the emulated engine cannot find a
# "whole word" boundary against a leading/trailing space and Execute() will
# simply return $false (no match, no replace) in that case.

# --- Paragraph: "Stack manipulation operations ..." ---
$r = $paraStack.Range
$r.Find.ClearFormatting()
$null = $r.Find.Execute("syntactic sugar", $true, $false, $false, $false, $false, $true, 1, $false, "core", 2)

$r = $paraStack.Range
$r.Find.ClearFormatting()
$null = $r.Find.Execute(" level, because the operations will be defined using operations on basic data types implemented at the core level.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " level, because the operations will be defined as a part of a command", 2)

# --- Paragraph: "Procedures/functions with arguments ..." ---
$r = $paraProc.Range
$r.Find.ClearFormatting()
$null = $r.Find.Execute("syntactic sugar", $true, $false, $false, $false, $false, $true, 1, $false, "core", 2)

$r = $paraProc.Range
$r.Find.ClearFormatting()
$null = $r.Find.Execute(" level, because", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " level, because they will be a part of the basic command implementation.", 2)

# --- Paragraph: "Static type system ..." ---
$r = $paraStatic.Range
$r.Find.ClearFormatting()
$null = $r.Find.Execute("syntactic sugar", $true, $false, $false, $false, $false, $true, 1, $false, "core", 2)

$r = $paraStatic.Range
$r.Find.ClearFormatting()
$null = $r.Find.Execute(" level, because the type checking will be built upon the basic data types.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " level, because the type checking will be built with the basic data types.", 2)
